$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the header text above the second data block (row 16) ---
# Old text ("thread count=4") moves down to the new block at row 24,
# and row 16 gets the new text ("thread count=1").
$ws.Range("A24").Value = "Env: Multiple threads for reading (thread count=4) and single thread for writes. Version 0.0.3"
$ws.Range("A16").Value = "Env: Multiple threads for reading (thread count=1) and single thread for writes. Version 0.0.3"

# --- New timing data rows 25-28 (columns A/B with formula in C) ---
$timeFmt = "h:mm:ss"

$ws.Range("A25").Value = 0.49806712962962968
$ws.Range("A25").NumberFormat = $timeFmt
$ws.Range("B25").Value = 0.49835648148148143
$ws.Range("B25").NumberFormat = $timeFmt

$ws.Range("A26").Value = 0.50166666666666659
$ws.Range("A26").NumberFormat = $timeFmt
$ws.Range("B26").Value = 0.5019675925925926
$ws.Range("B26").NumberFormat = $timeFmt

$ws.Range("A27").Value = 0.50302083333333336
$ws.Range("A27").NumberFormat = $timeFmt
$ws.Range("B27").Value = 0.50331018518518522
$ws.Range("B27").NumberFormat = $timeFmt

$ws.Range("A28").Value = 0.50428240740740737
$ws.Range("A28").NumberFormat = $timeFmt
$ws.Range("B28").Value = 0.50457175925925923
$ws.Range("B28").NumberFormat = $timeFmt

# Shared formula for C25:C28
$ws.Range("C25:C28").Formula = "=B25-A25"
$ws.Range("C25:C28").NumberFormat = $timeFmt

# Empty, but time-formatted, placeholder cells in G/H for rows 25-28
$ws.Range("G25:H28").NumberFormat = $timeFmt

# Shared formula for I25:I28
$ws.Range("I25:I28").Formula = "=H25-G25"
$ws.Range("I25:I28").NumberFormat = $timeFmt

# --- Averages row 29 ---
$ws.Range("C29").Formula = "=AVERAGE(C25:C28)"
$ws.Range("C29").NumberFormat = $timeFmt
$ws.Range("I29").Formula = "=AVERAGE(I25:I28)"
$ws.Range("I29").NumberFormat = $timeFmt

# --- Update selection to match the author's final cursor position ---
$ws.Range("B29").Select() | Out-Null
